$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B/C shift right to C/D).
$ws.Columns("B:B").Insert() | Out-Null

# New column B should be as wide as column A (both carry the wide
# "75.81640625" width used for the long query columns).
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth()

# Header row: B1 becomes the new "StatQuery" label.
$ws.Range("B1").Value = "StatQuery"

# Row 2: B2 gets the companion OPTIONAL MATCH stats query, wrapped the
# same way as A2 (style index 1 / "Normal 2" wrap-text style).
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Gall bladder carcinoma (adenocarcinoma)'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# Match the author's final selection/active cell.
$ws.Range("A2").Select() | Out-Null
